$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: keep value 1, change style to right-aligned (style index 6)
$ws.Range("A2").Value = 1
$ws.Range("A2").HorizontalAlignment = -4152

# A3: value changes from 9999 to -1, style becomes right-aligned (style index 6)
$ws.Range("A3").Value = -1
$ws.Range("A3").HorizontalAlignment = -4152

# A5: was empty, now holds 9999 with right-aligned style (style index 6)
$ws.Range("A5").Value = 9999
$ws.Range("A5").HorizontalAlignment = -4152

# A6: was empty, now holds text "null" with right-aligned style (style index 6)
$ws.Range("A6").Value = "null"
$ws.Range("A6").HorizontalAlignment = -4152

# Update the selected/active cell to C7
$ws.Range("C7").Select()

# Set page setup: letter/A4 sized paper (9 = A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
